$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $rng = $ws.Range($addr)
    $style = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $style
}

Set-CellText 'D2' '30.756.81'
Set-CellText 'E2' '  +0.96%  '
Set-CellText 'D3' '2.150.40'
Set-CellText 'E3' '  +2.12%  '
Set-CellText 'E4' '  +0.52%  '
Set-CellText 'D5' '352.26'
Set-CellText 'E5' '  +5.47%  '
Set-CellText 'E6' '  +0.40%  '
Set-CellText 'D7' '0.5287'
Set-CellText 'E7' '  +1.24%  '
Set-CellText 'D8' '0.4573'
Set-CellText 'E8' '  +1.02%  '
Set-CellText 'D9' '54.28'
Set-CellText 'E9' '  +2.09%  '
Set-CellText 'D10' '0.09205'
Set-CellText 'E10' '  +3.30%  '
Set-CellText 'D11' '1.187'
Set-CellText 'E11' '  +0.77%  '
Set-CellText 'D12' '24.94'
Set-CellText 'E12' '  +3.51%  '
Set-CellText 'D13' '2.152.71'
Set-CellText 'E13' '  +2.41%  '
Set-CellText 'D14' '6.921'
Set-CellText 'D15' '8.179'
Set-CellText 'E15' '  +1.85%  '
Set-CellText 'D16' '102.42'
Set-CellText 'E16' '  +5.90%  '
Set-CellText 'D17' '0.00001181'
Set-CellText 'E17' '  +3.37%  '
Set-CellText 'D18' '1.010'
Set-CellText 'E18' '  +0.39%  '
Set-CellText 'D19' '0.06721'
Set-CellText 'E19' '  +1.04%  '
Set-CellText 'D20' '19.64'
Set-CellText 'E20' '  +2.20%  '
Set-CellText 'E21' '  +0.36%  '
Set-CellText 'D22' '6.397'
Set-CellText 'E22' '  +0.89%  '
Set-CellText 'D23' '30.832.42'
Set-CellText 'E23' '  +1.01%  '
Set-CellText 'E24' '  +3.72%  '
Set-CellText 'E25' '  +1.96%  '
Set-CellText 'D26' '2.396.76'
Set-CellText 'E26' '  +1.95%  '
Set-CellText 'D27' '22.72'
Set-CellText 'E27' '  +2.26%  '
Set-CellText 'D28' '2.620'
Set-CellText 'E28' '  +3.60%  '
Set-CellText 'D29' '165.00'
Set-CellText 'E29' '  +1.33%  '
Set-CellText 'E30' '  +2.91%  '
Set-CellText 'D31' '1.228'
Set-CellText 'E31' '  +1.71%  '
Set-CellText 'D32' '0.1085'
Set-CellText 'E32' '  +1.31%  '
Set-CellText 'D33' '1.681'
Set-CellText 'E33' '  +1.20%  '
Set-CellText 'D34' '6.416'
Set-CellText 'E34' '  -0.16%  '
Set-CellText 'D35' '4.012'
Set-CellText 'E35' '  +1.83%  '
Set-CellText 'D36' '6.180'
Set-CellText 'E36' '  +6.51%  '
Set-CellText 'D37' '10.39'
Set-CellText 'E37' '  -0.34%  '
Set-CellText 'D38' '0.02662'
Set-CellText 'E38' '  +2.60%  '
Set-CellText 'D39' '0.06935'
Set-CellText 'E39' '  +1.31%  '
Set-CellText 'D40' '0.2343'
Set-CellText 'E40' '  +1.88%  '
Set-CellText 'D41' '12.72'
Set-CellText 'E41' '  -0.07%  '
Set-CellText 'D42' '0.6975'
Set-CellText 'E42' '  +1.46%  '
Set-CellText 'E43' '  +1.98%  '
Set-CellText 'D44' '14.94'
Set-CellText 'E44' '  +6.78%  '
Set-CellText 'D45' '2.367'
Set-CellText 'E45' '  +2.40%  '
Set-CellText 'D46' '0.6497'
Set-CellText 'E46' '  +2.06%  '
Set-CellText 'B47' 'PancakeSwap'
Set-CellText 'C47' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText 'D47' '3.758'
Set-CellText 'E47' '  +2.62%  '
Set-CellText 'B48' 'BabyDogeCoin'
Set-CellText 'C48' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText 'D48' '0.00000000369'
Set-CellText 'E48' '  +5.31%  '
Set-CellText 'D49' '1.260'
Set-CellText 'E49' '  +0.94%  '
Set-CellText 'D50' '83.62'
Set-CellText 'E50' '  +0.15%  '
Set-CellText 'D51' '0.07338'
Set-CellText 'E51' '  +2.71%  '

$wb.Save()
